# Fruta / hortaliza, semanal
# Insert 3 new weekly price rows for "Macroferia Regional de Talca - Limón"
# right before the existing row 768, shifting the remaining rows (old 768-817)
# down to 771-820, and fill in the new rows' data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three blank rows at 768..770 (existing content shifts down).
$ws.Rows("768:770").Insert()

# Row 768
$ws.Range("A768").Value = 5
$ws.Range("B768").Value = "Macroferia Regional de Talca"
$ws.Range("C768").Value = "Maule"
$ws.Range("D768").Value = 44516
$ws.Range("E768").Value = 7
$ws.Range("F768").Value = "Fruta"
$ws.Range("G768").Value = 100102
$ws.Range("H768").Value = "Cítricos"
$ws.Range("I768").Value = 100102003
$ws.Range("J768").Value = "Limón"
$ws.Range("K768").Value = "Sin especificar"
$ws.Range("L768").Value = "1a amarillo"
$ws.Range("M768").Value = 350
$ws.Range("N768").Value = 6000
$ws.Range("O768").Value = 6000
$ws.Range("P768").Value = 6000
$ws.Range("Q768").Value = "`$/malla 14 kilos"
$ws.Range("R768").Value = "Provincia de Quillota"
$ws.Range("S768").Value = 429
$ws.Range("T768").Value = 14

# Row 769
$ws.Range("A769").Value = 5
$ws.Range("B769").Value = "Macroferia Regional de Talca"
$ws.Range("C769").Value = "Maule"
$ws.Range("D769").Value = 44516
$ws.Range("E769").Value = 7
$ws.Range("F769").Value = "Fruta"
$ws.Range("G769").Value = 100102
$ws.Range("H769").Value = "Cítricos"
$ws.Range("I769").Value = 100102003
$ws.Range("J769").Value = "Limón"
$ws.Range("K769").Value = "Sin especificar"
$ws.Range("L769").Value = "1a amarillo"
$ws.Range("M769").Value = 720
$ws.Range("N769").Value = 5500
$ws.Range("O769").Value = 6000
$ws.Range("P769").Value = 5861
$ws.Range("Q769").Value = "`$/malla 14 kilos"
$ws.Range("R769").Value = "Región de O'Higgins"
$ws.Range("S769").Value = 419
$ws.Range("T769").Value = 14

# Row 770
$ws.Range("A770").Value = 5
$ws.Range("B770").Value = "Macroferia Regional de Talca"
$ws.Range("C770").Value = "Maule"
$ws.Range("D770").Value = 44516
$ws.Range("E770").Value = 7
$ws.Range("F770").Value = "Fruta"
$ws.Range("G770").Value = 100102
$ws.Range("H770").Value = "Cítricos"
$ws.Range("I770").Value = 100102003
$ws.Range("J770").Value = "Limón"
$ws.Range("K770").Value = "Sin especificar"
$ws.Range("L770").Value = "3a amarillo"
$ws.Range("M770").Value = 140
$ws.Range("N770").Value = 3500
$ws.Range("O770").Value = 3500
$ws.Range("P770").Value = 3500
$ws.Range("Q770").Value = "`$/malla 14 kilos"
$ws.Range("R770").Value = "Región de O'Higgins"
$ws.Range("S770").Value = 250
$ws.Range("T770").Value = 14
